# Append: 2026-02-13 12:54 JST
# A new job listing (score 135) is inserted, in descending-score order, between
# the existing rows that scored 218 and 55 - i.e. at row 7 - pushing the
# previous rows 7-11 down to rows 8-12. Every row's "取得日時" (fetch
# timestamp) in column A is refreshed to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-02-13 12:54:34"

# Shift rows 7:11 down to 8:12, carrying their formatting/hyperlinks with them,
# and open up row 7 for the newly scraped listing.
$ws.Rows("7:7").Insert()

# New listing slotted into row 7 (score 135 sits between 218 and 55).
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5473940"
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = "◆開発 ◇業務改善"

# Refresh the fetch timestamp on every other existing row (2-6 unchanged
# content, 8-12 are the old rows 7-11 shifted down by the insert).
$ws.Range("A2").Value = $timestamp
$ws.Range("A3").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A8").Value = $timestamp
$ws.Range("A9").Value = $timestamp
$ws.Range("A10").Value = $timestamp
$ws.Range("A11").Value = $timestamp
$ws.Range("A12").Value = $timestamp

# Row 12 (previously row 11) needs its own hyperlink entry added, since the
# insert only carried the existing 10 hyperlinks down with their rows - it
# did not fabricate an 11th one for the newly exposed last row.
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5490905")
$ws.Range("F12").Style = "Hyperlink"
